$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("28")

$values = @{
  2 = @(1.836570346695847, 1.96036690553072)
  3 = @(2.714261191003767, 3.877085295812139)
  4 = @(5.311891470459845, 5.965292690314715)
  5 = @(6.15422784167213, 8.016623461272442)
  6 = @(9.85293134482966, 9.992153466566959)
  7 = @(10.4712054436148, 12.24408646052714)
  8 = @(16.02307573052934, 14.14681200848843)
  9 = @(17.34701051379703, 16.1421766770639)
  10 = @(17.63930664055301, 18.13190342498724)
  11 = @(20.40644316777218, 20.25672230162973)
  12 = @(22.50801094171332, 22.63859122061831)
  13 = @(28.19635979995576, 24.51205309412098)
  14 = @(30.0158194680641, 26.28133380295972)
  15 = @(30.99487649024563, 28.4091613865959)
  16 = @(32.76801620776994, 30.43487046935713)
  17 = @(35.55660041797303, 32.62877408201998)
  18 = @(37.818976932774, 34.6255795201146)
  19 = @(38.75320938252725, 36.37284374515027)
  20 = @(39.57724138340481, 38.44066646180399)
  21 = @(39.83469695094952, 40.65711678213948)
  22 = @(39.93376651393991, 42.56394554571474)
  23 = @(40.30138110509802, 44.61665579937293)
  24 = @(42.05886214667142, 46.36668027265851)
  25 = @(43.81387140733236, 48.33727649334537)
  26 = @(48.50411245895745, 50.37800334206129)
  27 = @(54.30951254556948, 52.14866232195176)
  28 = @(57.03939531926715, 54.08270892561379)
  29 = @(57.98630815488951, 56.00385695311693)
  30 = @(60.99500471085199, 57.91294066701746)
  31 = @(63.40974947496235, 59.8996195904618)
  32 = @(64.88555087930017, 62.01882427224799)
  33 = @(65.28759822568567, 63.92312366625588)
  34 = @(66.45436210475127, 65.87404079025841)
  35 = @(69.25998582639978, 68.37567231043249)
  36 = @(74.34239320165251, 70.22361291175963)
  37 = @(76.75561101239953, 72.31360563528949)
  38 = @(80.04769583365753, 74.32509059628343)
  39 = @(81.41304005115973, 76.2335061477208)
  40 = @(82.26571891056187, 78.03996516697632)
  41 = @(82.51746998987669, 79.79290720682204)
  42 = @(83.60169195955582, 81.91606523947326)
  43 = @(85.68747994543715, 83.69815773714902)
  44 = @(86.51911434575526, 85.749628110865)
  45 = @(91.34926719130266, 87.56096242341324)
  46 = @(93.4958875262132, 89.83542655597611)
  47 = @(94.42755217681183, 91.81763896012852)
  48 = @(95.72448597396846, 93.82590699549186)
  49 = @(99.93677721042408, 96.00151164036893)
}

foreach ($r in $values.Keys) {
  $pair = $values[$r]
  $ws.Cells.Item([int]$r, 2).Value = $pair[0]
  $ws.Cells.Item([int]$r, 3).Value = $pair[1]
}
